$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / percentage cells: Coin names, Links, Volume(1h) percentages.
# These never look like pure numbers, so a normal .Value assignment keeps them as text.
$textUpdates = @{
    "B10" = "BinanceUSD"
    "C10" = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
    "B11" = "OKB"
    "C11" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "B29" = "Filecoin"
    "C29" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "B30" = "BitcoinCash"
    "C30" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "E3" = "  +0.11%  "
    "E4" = "  +0.26%  "
    "E5" = "  +0.09%  "
    "E7" = "  -0.19%  "
    "E9" = "  -2.09%  "
    "E10" = "  +0.39%  "
    "E11" = "  -2.79%  "
    "E12" = "  +0.98%  "
    "E13" = "  -0.48%  "
    "E14" = "  +0.11%  "
    "E15" = "  +4.68%  "
    "E16" = "  -0.24%  "
    "E17" = "  -0.10%  "
    "E18" = "  -0.74%  "
    "E19" = "  -0.30%  "
    "E20" = "  -0.22%  "
    "E21" = "  +1.37%  "
    "E22" = "  +0.17%  "
    "E23" = "  +3.15%  "
    "E24" = "  +0.61%  "
    "E25" = "  +6.79%  "
    "E26" = "  +1.15%  "
    "E27" = "  +1.44%  "
    "E28" = "  +1.22%  "
    "E29" = "  +4.50%  "
    "E30" = "  +1.90%  "
    "E31" = "  -1.16%  "
    "E32" = "  +2.51%  "
    "E33" = "  -2.63%  "
    "E34" = "  -1.87%  "
    "E35" = "  +0.04%  "
    "E36" = "  +1.39%  "
    "E37" = "  -0.59%  "
    "E38" = "  -2.81%  "
    "E39" = "  +1.69%  "
    "E40" = "  -1.58%  "
    "E42" = "  -0.93%  "
    "E43" = "  +1.59%  "
    "E44" = "  -1.59%  "
    "E45" = "  +2.09%  "
    "E46" = "  +0.07%  "
    "E47" = "  +0.30%  "
    "E48" = "  -0.85%  "
    "E49" = "  -4.64%  "
    "E50" = "  +2.74%  "
    "E51" = "  -0.38%  "
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Price cells: many look like plain numbers (e.g. "1.003"), which Excel would
# otherwise auto-convert to a numeric value. Force text format first, assign the
# literal string, then restore the default 'Normal' style so no formatting residue
# is left behind on the cell.
$priceUpdates = @{
    "D2" = "24.798.61"
    "D3" = "1.702.66"
    "D4" = "1.003"
    "D5" = "316.39"
    "D7" = "0.3938"
    "D8" = "0.4045"
    "D9" = "1.510"
    "D10" = "1.004"
    "D11" = "53.62"
    "D12" = "0.08908"
    "D13" = "7.275"
    "D14" = "23.46"
    "D15" = "8.002"
    "D16" = "0.00001331"
    "D17" = "1.699.43"
    "D18" = "100.21"
    "D19" = "0.07044"
    "D20" = "19.68"
    "D21" = "7.033"
    "D23" = "14.64"
    "D24" = "24.781.76"
    "D25" = "3.199"
    "D27" = "22.77"
    "D28" = "162.05"
    "D29" = "8.125"
    "D30" = "136.69"
    "D31" = "5.175"
    "D32" = "0.08804"
    "D33" = "1.083"
    "D34" = "7.302"
    "D35" = "11.21"
    "D36" = "1.979"
    "D37" = "0.2747"
    "D38" = "14.41"
    "D39" = "0.09213"
    "D40" = "0.02753"
    "D41" = "1.463"
    "D42" = "0.7700"
    "D43" = "15.86"
    "D44" = "0.7181"
    "D45" = "2.573"
    "D46" = "4.212"
    "D48" = "140.69"
    "D49" = "1.313"
    "D50" = "90.91"
    "D51" = "0.08006"
}
foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
    $cell.Style = "Normal"
}
